$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (dates replacing years in column A, refreshed figures in B:E)
$data = @(
    @(40816, 239.8890918198565,  220.24823461320381,  14.447074297902599,   474.5844007309629),
    @(41182, 195.0309780261266,  304.2693837266757,   8.5172113222631207,   507.81757307506541),
    @(41547, -14.699067618251346,234.35071452721741,  -28.599272865945551,  191.05237404302051),
    @(41912, -412.74994654988103,224.04280647608536,  159.7779401013828,    -28.929199972412874),
    @(42277, 432.70105978629772, 266.27757525310625,  591.21038404956835,   1290.1890190889724),
    @(42643, 137.33185378059758, 274.67789034068295,  1237.2773169587701,   1649.2870610800505),
    @(43008, -874.2294138262082, 196.33205755767355,  -585.40044837138896,  -1263.2978046399235),
    @(43373, 547.71937948598304, 199.21655750978803,  -122.83383262001995,  624.10210437575108),
    @(43738, 182.2331408862365,  182.742385376627,    435.21307259997178,   800.18859886283531),
    @(44104, 705.47815715598517, 117.77661656406961,  1314.6788939662877,   2137.9336676863422),
    @(44469, -21.990848899680397,60.449932908935246,  50.179420690060859,   88.638504699315718),
    @(44834, -211.64866738028053,169.6851080692422,   -73.025590623547103,  -114.98914993458543),
    @(45199, -88.66,             219.96000000000004,  -97.26,               34.040000000000035)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}

# Column A now holds dates rather than plain years -> format as a date
# (numFmtId 14, Excel's built-in short-date format)
$ws.Range("A2:A14").NumberFormat = "mm-dd-yy"
